# feat: add 2022-Q4 data
#
# - "总计": push the existing "2020-Q4" summary row down to row 3 and add a
#   new row 2 summarising "2022-Q4".
# - Insert a new "2022-Q4" sheet (cloned from "总计", matching the page
#   margins/header style the author reused) between "总计" and "2020-Q4",
#   and fill it in with the fund holdings table.
# - Leave "2020-Q4" (now the 3rd tab) untouched and re-select it as the
#   active tab.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q4_2020 = $wb.Worksheets.Item("2020-Q4")

# ---------------------------------------------------------------------
# 1. "总计": shift the old "2020-Q4" row down to row 3, add new row 2 for
#    "2022-Q4".
# ---------------------------------------------------------------------
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2020-Q4"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.16

$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.01

# ---------------------------------------------------------------------
# 2. Clone "总计" into a new sheet positioned right before "2020-Q4" (this
#    mirrors the original author's workflow: same page margins/header
#    style as "总计") and rename it to "2022-Q4".
# ---------------------------------------------------------------------
$totalSheet.Copy($q4_2020)
$newSheet = $wb.Worksheets.Item("总计 (2)")
$newSheet.Name = "2022-Q4"

# Clear the cloned "总计" data out of the new sheet before filling it in,
# but keep the header/index-column styling (s=2) that came along with it.
$newSheet.Range("A1:D2").ClearContents()

$newSheet.Range("B1").Copy()
$newSheet.Range("E1:H1").PasteSpecial(-4122)

$newSheet.Range("A2").Copy()
$newSheet.Range("A3").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0

$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "013920"
$newSheet.Range("B2").ClearFormats()

$newSheet.Range("C2").Value = "兴华创新医疗6个月持有混合A"

$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.22"
$newSheet.Range("D2").ClearFormats()

$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "92.79"
$newSheet.Range("E2").ClearFormats()

$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "4.51"
$newSheet.Range("F2").ClearFormats()

$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0099"
$newSheet.Range("G2").ClearFormats()

$newSheet.Range("H2").Value = 4

$newSheet.Range("A3").Value = 1

$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "013921"
$newSheet.Range("B3").ClearFormats()

$newSheet.Range("C3").Value = "兴华创新医疗6个月持有混合C"

$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "0.06"
$newSheet.Range("D3").ClearFormats()

$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "92.79"
$newSheet.Range("E3").ClearFormats()

$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "4.51"
$newSheet.Range("F3").ClearFormats()

$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0027"
$newSheet.Range("G3").ClearFormats()

$newSheet.Range("H3").Value = 4

# ---------------------------------------------------------------------
# 3. Restore "2020-Q4" as the active tab (matches the unchanged sheet3.xml
#    in the diff, which still carries tabSelected="1").
# ---------------------------------------------------------------------
$active2020 = $wb.Worksheets.Item("2020-Q4")
$active2020.Activate()
